# Update OpenData_Slovakia_Covid_DailyStats sheet1: po 19. 07. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 393-498: set AgTests (F) / AgPosit (G) values (new columns for older rows,
# corrected/updated values for most recent rows already having F/G).
$fgData = @(
    @(393, 309055, 1242),
    @(394, 166076, 623),
    @(395, 753986, 1956),
    @(396, 166408, 549),
    @(397, 107953, 639),
    @(398, 300481, 1477),
    @(399, 200568, 967),
    @(400, 149862, 803),
    @(401, 272927, 937),
    @(402, 723063, 1392),
    @(403, 353970, 734),
    @(404, 224848, 914),
    @(405, 174931, 693),
    @(406, 171548, 680),
    @(407, 158330, 674),
    @(408, 305160, 835),
    @(409, 708940, 1007),
    @(410, 364822, 635),
    @(411, 225429, 828),
    @(412, 176630, 647),
    @(413, 149804, 658),
    @(414, 149258, 563),
    @(415, 308617, 695),
    @(416, 672323, 931),
    @(417, 343404, 589),
    @(418, 202443, 701),
    @(419, 149708, 510),
    @(420, 138984, 501),
    @(421, 153274, 534),
    @(422, 298724, 646),
    @(423, 439862, 636),
    @(424, 266318, 497),
    @(425, 138148, 549),
    @(426, 107381, 382),
    @(427, 90523, 370),
    @(428, 102463, 389),
    @(429, 178138, 458),
    @(430, 175311, 272),
    @(431, 171439, 400),
    @(432, 122694, 425),
    @(433, 86929, 269),
    @(434, 79053, 279),
    @(435, 83440, 266),
    @(436, 145398, 353),
    @(437, 167330, 272),
    @(438, 121829, 250),
    @(439, 89308, 318),
    @(440, 73695, 226),
    @(441, 68285, 202),
    @(442, 70542, 172),
    @(443, 106930, 208),
    @(444, 104111, 192),
    @(445, 84526, 189),
    @(446, 86698, 264),
    @(447, 66963, 190),
    @(448, 61456, 140),
    @(449, 59895, 155),
    @(450, 91439, 170),
    @(451, 86521, 120),
    @(452, 74575, 126),
    @(453, 70159, 210),
    @(454, 52585, 133),
    @(455, 50661, 119),
    @(456, 50351, 133),
    @(457, 78972, 135),
    @(458, 70686, 74),
    @(459, 59749, 87),
    @(460, 58501, 148),
    @(461, 45203, 94),
    @(462, 43878, 53),
    @(463, 46674, 72),
    @(464, 73302, 84),
    @(465, 61288, 57),
    @(466, 51014, 57),
    @(467, 52210, 78),
    @(468, 41584, 47),
    @(469, 40967, 39),
    @(470, 43416, 43),
    @(471, 66533, 54),
    @(472, 51185, 22),
    @(473, 40299, 42),
    @(474, 45037, 59),
    @(475, 36119, 29),
    @(476, 37129, 30),
    @(477, 37729, 38),
    @(478, 54303, 33),
    @(479, 41833, 32),
    @(480, 33837, 22),
    @(481, 43561, 36),
    @(482, 36110, 23),
    @(483, 65586, 36),
    @(484, 8247, 11),
    @(485, 13872, 14),
    @(486, 8892, 8),
    @(487, 6821, 9),
    @(488, 6362, 9),
    @(489, 12534, 11),
    @(490, 10658, 15),
    @(491, 9852, 11),
    @(492, 13791, 17),
    @(493, 8110, 7),
    @(494, 6148, 6),
    @(495, 10207, 13),
    @(496, 7909, 14),
    @(497, 7305, 9),
    @(498, 8732, 9)
)

foreach ($row in $fgData) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
}

# Rows 499-501: brand-new daily rows appended at the bottom of the table.
$newRows = @(
    @(499, 44393, 392071, 9384, 37, 12524, 9583, 10),
    @(500, 44394, 392100, 6275, 29, 12524, 6034, 5),
    @(501, 44395, 392104, 1168, 4, 12527, 4544, 6)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
